$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "HK_R_acc_G"
$ws.Range("A2").Value = 79.294478527607353
$ws.Range("A3").Value = 78.834355828220865
$ws.Range("A4").Value = 75.920245398773005
$ws.Range("A5").Value = 76.99386503067484
$ws.Range("A6").Value = 76.917177914110425
$ws.Range("A7").Value = 79.064417177914109
$ws.Range("A8").Value = 76.763803680981596
$ws.Range("A9").Value = 76.840490797546011
$ws.Range("A10").Value = 76.303680981595093
$ws.Range("A11").Value = 76.303680981595093
$ws.Range("A12").Value = 77.760736196319016
$ws.Range("A13").Value = 77.607361963190186
$ws.Range("A14").Value = 77.607361963190186
$ws.Range("A15").Value = 77.607361963190186
$ws.Range("A16").Value = 77.530674846625772
$ws.Range("A17").Value = 76.457055214723923
$ws.Range("A18").Value = 78.604294478527606
$ws.Range("A19").Value = 78.604294478527606
$ws.Range("A20").Value = 80.521472392638032
$ws.Range("A21").Value = 77.147239263803684
$ws.Range("A22").Value = 80.598159509202446
$ws.Range("A23").Value = 78.144171779141104
$ws.Range("A24").Value = 77.914110429447859
$ws.Range("A25").Value = 77.914110429447859
$ws.Range("A26").Value = 76.457055214723923
$ws.Range("A27").Value = 76.457055214723923
$ws.Range("A28").Value = 75.766871165644162
$ws.Range("A29").Value = 75.766871165644162
$ws.Range("A30").Value = 75.843558282208591
$ws.Range("A31").Value = 75.99693251533742
$ws.Range("A32").Value = 75.613496932515332
$ws.Range("A33").Value = 76.150306748466249
$ws.Range("A34").Value = 75.766871165644162
$ws.Range("A35").Value = 75.843558282208591
$ws.Range("A36").Value = 75.920245398773005
$ws.Range("A37").Value = 82.898773006134974
$ws.Range("A38").Value = 75.920245398773005
$ws.Range("A39").Value = 75.766871165644162
$ws.Range("A40").Value = 77.914110429447859
$ws.Range("A41").Value = 75.690184049079761
$ws.Range("A42").Value = 76.457055214723923
$ws.Range("A43").Value = 75.766871165644162
$ws.Range("A44").Value = 75.766871165644162
$ws.Range("A45").Value = 76.303680981595093
$ws.Range("A46").Value = 75.690184049079761
$ws.Range("A47").Value = 75.306748466257673
$ws.Range("A48").Value = 76.073619631901849
$ws.Range("A49").Value = 75.383435582822088
